$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bring format (number format / font) for the new row's B and C cells
# from the row above, reusing the existing style definitions instead of
# creating new ones (Insert with a copied clipboard applies formats only
# to the not-yet-used destination cells).
$ws.Range("B26:C26").Copy()
$ws.Range("B27:C27").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# Fill in the new row of data (Lote 69553031 - Bruno, bonus of 5, on 27/06/2025)
$ws.Range("A27").Value = 25
$ws.Range("B27").Value = 45835
$ws.Range("C27").Value = 69553031
$ws.Range("D27").Value = 5
$ws.Range("E27").Value = "Bruno"
$ws.Range("F27").Value = "Bruno fez sozinho"

$ws.Range("D27").Select()
